$wb = $excel.ActiveWorkbook

# zh-cn sheet: Correspond Handoff Datetime (D21:D22) and Correspond Handback DateTime (G21:G22)
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D21:D22").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZhCn.Range("D21:D22").Value = "2016-02-29 04:36:20"
$wsZhCn.Range("G21:G22").Value = "2016-02-29 04:37:13"

# de-de sheet: Correspond Handoff Datetime (D21:D22) and Correspond Handback DateTime (G21:G22)
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D21:D22").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDeDe.Range("D21:D22").Value = "2016-02-29 04:36:36"
$wsDeDe.Range("G21:G22").Value = "2016-02-29 04:37:35"
